$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 22158.8
$ws.Cells.Item(58, 10).Value = 25135
$ws.Cells.Item(58, 12).Value = 75405
$ws.Cells.Item(58, 14).Value = -75705
$ws.Cells.Item(135, 8).Value = 620.9792
$ws.Cells.Item(135, 9).Value = 519.0238000000001
$ws.Cells.Item(135, 10).Value = 1334.6666
$ws.Cells.Item(135, 11).Value = 4671.2142
$ws.Cells.Item(135, 12).Value = 12011.9994
$ws.Cells.Item(135, 13).Value = -2136.2142
$ws.Cells.Item(135, 14).Value = -17081.9994
$ws.Cells.Item(137, 8).Value = 3826.2593
$ws.Cells.Item(137, 9).Value = 4199.263
$ws.Cells.Item(137, 10).Value = 2940.375
$ws.Cells.Item(137, 11).Value = 12597.789
$ws.Cells.Item(137, 12).Value = 8821.125
$ws.Cells.Item(137, 13).Value = -10047.789
$ws.Cells.Item(137, 14).Value = -13921.125
$ws.Cells.Item(138, 8).Value = 4302.9727
$ws.Cells.Item(138, 9).Value = 1498.8096
$ws.Cells.Item(138, 10).Value = 5435.423
$ws.Cells.Item(138, 11).Value = 4496.4288
$ws.Cells.Item(138, 12).Value = 16306.269
$ws.Cells.Item(138, 13).Value = 643.5712000000003
$ws.Cells.Item(138, 14).Value = -26586.269
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3186.3225
$ws.Cells.Item(61, 9).Value = 1141
$ws.Cells.Item(61, 10).Value = 4478.1055
$ws.Cells.Item(61, 11).Value = 1141
$ws.Cells.Item(61, 12).Value = 4478.1055
$ws.Cells.Item(61, 13).Value = -929
$ws.Cells.Item(61, 14).Value = -4902.1055
$ws.Cells.Item(74, 8).Value = 663.8148
$ws.Cells.Item(74, 10).Value = 928.5
$ws.Cells.Item(74, 12).Value = 928.5
$ws.Cells.Item(74, 14).Value = -2676.5
$ws.Cells.Item(77, 8).Value = 663.8148
$ws.Cells.Item(77, 10).Value = 928.5
$ws.Cells.Item(77, 12).Value = 4642.5
$ws.Cells.Item(77, 14).Value = -13378.5
$ws.Cells.Item(132, 8).Value = 16951078
$ws.Cells.Item(132, 9).Value = 21277804
$ws.Cells.Item(132, 10).Value = 4735.5
$ws.Cells.Item(132, 11).Value = 63833412
$ws.Cells.Item(132, 12).Value = 14206.5
$ws.Cells.Item(132, 13).Value = -63830882
$ws.Cells.Item(132, 14).Value = -19266.5
$ws.Cells.Item(136, 8).Value = 3186.3225
$ws.Cells.Item(136, 9).Value = 1141
$ws.Cells.Item(136, 10).Value = 4478.1055
$ws.Cells.Item(136, 11).Value = 3423
$ws.Cells.Item(136, 12).Value = 13434.3165
$ws.Cells.Item(136, 13).Value = -873
$ws.Cells.Item(136, 14).Value = -18534.3165
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2302.5151
$ws.Cells.Item(134, 9).Value = 1835.7826
$ws.Cells.Item(134, 10).Value = 3376
$ws.Cells.Item(134, 11).Value = 5507.3478
$ws.Cells.Item(134, 12).Value = 10128
$ws.Cells.Item(134, 13).Value = -2972.3478
$ws.Cells.Item(134, 14).Value = -15198
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4534.8335
$ws.Cells.Item(31, 9).Value = 2582.9092
$ws.Cells.Item(31, 10).Value = 6186.4614
$ws.Cells.Item(31, 11).Value = 2582.9092
$ws.Cells.Item(31, 12).Value = 6186.4614
$ws.Cells.Item(31, 13).Value = -2287.9092
$ws.Cells.Item(31, 14).Value = -6776.4614
$ws.Cells.Item(34, 8).Value = 4534.8335
$ws.Cells.Item(34, 9).Value = 2582.9092
$ws.Cells.Item(34, 10).Value = 6186.4614
$ws.Cells.Item(34, 11).Value = 2582.9092
$ws.Cells.Item(34, 12).Value = 6186.4614
$ws.Cells.Item(34, 13).Value = -2380.9092
$ws.Cells.Item(34, 14).Value = -6590.4614
$ws.Cells.Item(99, 8).Value = 2031
$ws.Cells.Item(99, 9).Value = 1199.8
$ws.Cells.Item(99, 10).Value = 2550.5
$ws.Cells.Item(99, 11).Value = 1199.8
$ws.Cells.Item(99, 12).Value = 2550.5
$ws.Cells.Item(99, 13).Value = 298.2
$ws.Cells.Item(99, 14).Value = -5546.5
$ws.Cells.Item(107, 8).Value = 1511.2858
$ws.Cells.Item(107, 9).Value = 665
$ws.Cells.Item(107, 10).Value = 2639.6667
$ws.Cells.Item(107, 11).Value = 665
$ws.Cells.Item(107, 12).Value = 2639.6667
$ws.Cells.Item(107, 13).Value = 1255
$ws.Cells.Item(107, 14).Value = -6479.6667
$ws.Cells.Item(126, 8).Value = 2031
$ws.Cells.Item(126, 9).Value = 1199.8
$ws.Cells.Item(126, 10).Value = 2550.5
$ws.Cells.Item(126, 11).Value = 3599.4
$ws.Cells.Item(126, 12).Value = 7651.5
$ws.Cells.Item(126, 13).Value = -1129.4
$ws.Cells.Item(126, 14).Value = -12591.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 639.6
$ws.Cells.Item(17, 9).Value = 99.5
$ws.Cells.Item(17, 10).Value = 999.6667
$ws.Cells.Item(17, 11).Value = 298.5
$ws.Cells.Item(17, 12).Value = 2999.0001
$ws.Cells.Item(17, 13).Value = -129.5
$ws.Cells.Item(17, 14).Value = -3337.0001
$ws.Cells.Item(19, 8).Value = 3000
$ws.Cells.Item(19, 9).Value = 3000
$ws.Cells.Item(19, 10).Value = 3000
$ws.Cells.Item(19, 11).Value = 9000
$ws.Cells.Item(19, 12).Value = 9000
$ws.Cells.Item(22, 8).Value = 6338.769
$ws.Cells.Item(22, 9).Value = 1000.3333
$ws.Cells.Item(22, 10).Value = 7940.3
$ws.Cells.Item(22, 11).Value = 3000.9999
$ws.Cells.Item(22, 12).Value = 23820.9
$ws.Cells.Item(22, 14).Value = -24158.9
$ws.Cells.Item(27, 8).Value = 6338.769
$ws.Cells.Item(27, 9).Value = 1000.3333
$ws.Cells.Item(27, 10).Value = 7940.3
$ws.Cells.Item(27, 11).Value = 3000.9999
$ws.Cells.Item(27, 12).Value = 23820.9
$ws.Cells.Item(27, 14).Value = -24024.9
$ws.Cells.Item(32, 8).Value = 2199.7144
$ws.Cells.Item(32, 9).Value = 1332.6666
$ws.Cells.Item(32, 10).Value = 2850
$ws.Cells.Item(32, 11).Value = 3997.9998
$ws.Cells.Item(32, 12).Value = 8550
$ws.Cells.Item(32, 14).Value = -9116
$ws.Cells.Item(37, 8).Value = 37357
$ws.Cells.Item(37, 10).Value = 37357
$ws.Cells.Item(37, 12).Value = 112071
$ws.Cells.Item(37, 14).Value = -112295
$ws.Cells.Item(39, 8).Value = 1620.0526
$ws.Cells.Item(39, 10).Value = 1722.4117
$ws.Cells.Item(39, 12).Value = 5167.2351
$ws.Cells.Item(39, 14).Value = -5755.2351
$ws.Cells.Item(46, 8).Value = 2133.6365
$ws.Cells.Item(46, 9).Value = 2880
$ws.Cells.Item(46, 10).Value = 2059
$ws.Cells.Item(46, 11).Value = 8640
$ws.Cells.Item(46, 12).Value = 6177
$ws.Cells.Item(46, 13).Value = -8549
$ws.Cells.Item(46, 14).Value = -6359
$ws.Cells.Item(113, 8).Value = 1011173.44
$ws.Cells.Item(113, 9).Value = 18181818
$ws.Cells.Item(113, 10).Value = 1135.5294
$ws.Cells.Item(113, 11).Value = 54545454
$ws.Cells.Item(113, 12).Value = 3406.5882
$ws.Cells.Item(113, 13).Value = -54543284
$ws.Cells.Item(113, 14).Value = -7746.5882
$ws.Cells.Item(131, 8).Value = 1396.2963
$ws.Cells.Item(131, 10).Value = 1246.1111
$ws.Cells.Item(131, 12).Value = 3738.3333
$ws.Cells.Item(131, 14).Value = -13818.3333
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2845.375
$ws.Cells.Item(46, 9).Value = 350.25
$ws.Cells.Item(46, 10).Value = 5340.5
$ws.Cells.Item(46, 11).Value = 350.25
$ws.Cells.Item(46, 12).Value = 5340.5
$ws.Cells.Item(46, 13).Value = -162.25
$ws.Cells.Item(46, 14).Value = -5716.5
$ws.Cells.Item(132, 8).Value = 2478.5278
$ws.Cells.Item(132, 9).Value = 1126
$ws.Cells.Item(132, 10).Value = 5183.5835
$ws.Cells.Item(132, 11).Value = 3378
$ws.Cells.Item(132, 12).Value = 15550.7505
$ws.Cells.Item(132, 13).Value = -848
$ws.Cells.Item(132, 14).Value = -20610.7505
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(122, 8).Value = 359407.1
$ws.Cells.Item(122, 9).Value = 456368.88
$ws.Cells.Item(122, 10).Value = 3880.6667
$ws.Cells.Item(122, 11).Value = 1369106.64
$ws.Cells.Item(122, 12).Value = 11642.0001
$ws.Cells.Item(122, 13).Value = -1366656.64
$ws.Cells.Item(122, 14).Value = -16542.0001
$ws.Cells.Item(136, 8).Value = 821.32556
$ws.Cells.Item(136, 9).Value = 438.55884
$ws.Cells.Item(136, 10).Value = 2267.3333
$ws.Cells.Item(136, 11).Value = 1315.67652
$ws.Cells.Item(136, 12).Value = 6801.999899999999
$ws.Cells.Item(136, 13).Value = 1234.32348
$ws.Cells.Item(136, 14).Value = -11901.9999
$ws.Cells.Item(15, 14).ClearContents() | Out-Null
